$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 2
    14 = 2
    15 = 1
    16 = 2
    17 = 0
    18 = 2
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 1
    26 = 0
    27 = 0
    28 = 2
    29 = 1
    30 = 1
    31 = 0
    32 = 2
    33 = 2
    34 = 0
    35 = 3
    36 = 2
    37 = 1
    38 = 3
    40 = 2
    41 = 2
    42 = 2
    43 = 1
    44 = 1
    45 = 1
    46 = 2
    47 = 1
    48 = 0
    49 = 1
    51 = 0
    52 = 1
    53 = 2
    54 = 0
    55 = 1
    56 = 1
    57 = 2
    58 = 0
    59 = 2
    60 = 1
    61 = 2
    62 = 1
    64 = 3
    66 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
